$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price/Volume data cells so numeric-looking
# strings (e.g. "327.24") are not coerced into real numbers, matching
# the original inline-string cell types. Scoped to the data rows only so
# the header row (D1/E1) keeps its original General format.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.423.51"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "1.906.76"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "327.24"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "0.4679"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("D8").Value = "0.4084"
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").Value = "47.66"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").Value = "0.08001"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").Value = "1.005"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "22.29"
$ws.Range("E12").Value = "  +2.66%  "
$ws.Range("D13").Value = "1.919.88"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").Value = "5.928"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "7.113"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "89.02"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "0.06596"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").Value = "17.67"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D22").Value = "29.462.83"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").Value = "5.525"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("D25").Value = "2.207"
$ws.Range("D26").Value = "2.137.48"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "153.59"
$ws.Range("E27").Value = "  -2.55%  "
$ws.Range("D28").Value = "19.74"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").Value = "2.127"
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("D30").Value = "5.702"
$ws.Range("E30").Value = "  +5.58%  "
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("E32").Value = "  +9.39%  "
$ws.Range("D33").Value = "0.09487"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").Value = "3.572"
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("D36").Value = "5.376"
$ws.Range("E36").Value = "  +1.50%  "
$ws.Range("D37").Value = "0.02253"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").Value = "0.06075"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "8.341"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "1.170"
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("D41").Value = "0.5866"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D42").Value = "0.1835"
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "10.09"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "2.445"
$ws.Range("E44").Value = "  +4.46%  "
$ws.Range("D45").Value = "1.299"
$ws.Range("E45").Value = "  +3.21%  "
$ws.Range("D46").Value = "0.07738"
$ws.Range("E46").Value = "  +10.03%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.5532"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "12.11"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").Value = "1.927"
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("D50").Value = "113.20"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("E51").Value = "  +6.33%  "
